$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Through 2021-11-25"

$ws.Range("Q10").Value = 7
$ws.Range("R10").Value = 156
$ws.Range("S10").Value = 0.0429

$ws.Range("A13").Value = "November (through 11-25)"
$ws.Range("F13").Value = 55
$ws.Range("G13").Value = 0.1129
$ws.Range("I13").Value = 92
$ws.Range("J13").Value = 0.0213
$ws.Range("L13").Value = 43
$ws.Range("M13").Value = 0.14
$ws.Range("O13").Value = 39
$ws.Range("P13").Value = 0.1333
$ws.Range("R13").Value = 172
$ws.Range("S13").Value = 0.0444
$ws.Range("T13").Value = 4
$ws.Range("U13").Value = 167
$ws.Range("V13").Value = 0.0234

$ws.Range("F14").Value = 489
$ws.Range("G14").Value = 0.1077
$ws.Range("I14").Value = 741
$ws.Range("J14").Value = 0.0784
$ws.Range("L14").Value = 592
$ws.Range("M14").Value = 0.1098
$ws.Range("O14").Value = 473
$ws.Range("P14").Value = 0.1025
$ws.Range("Q14").Value = 61
$ws.Range("R14").Value = 1176
$ws.Range("S14").Value = 0.0493
$ws.Range("T14").Value = 96
$ws.Range("U14").Value = 1518
$ws.Range("V14").Value = 0.0595
